# Update the MPA test automation upload file: renumber asset/partner-asset
# identifiers on the "Data" sheet (rows 6-29).
#
# Value mapping applied to the relevant cells:
#   60000339 -> 60000347   (column K - ANLN1)
#   250      -> 257        (column L - ANLN2)
#   60000340 -> 60000348   (column N - PANL1)
#   251      -> 258        (column O - PANL2)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

for ($row = 6; $row -le 29; $row++) {
    foreach ($col in @("K", "L", "N", "O")) {
        $cell = $ws.Range("$col$row")
        $val = $cell.Value2

        if ($val -eq 60000339) {
            $cell.Value2 = 60000347
        }
        elseif ($val -eq 250) {
            $cell.Value2 = 257
        }
        elseif ($val -eq 60000340) {
            $cell.Value2 = 60000348
        }
        elseif ($val -eq 251) {
            $cell.Value2 = 258
        }
    }
}
